$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1: copy formatting (style) from H1 then set the text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows I2:J35
$data = @{
  2 = @(6, 6)
  3 = @(7, 7)
  4 = @(6, 6)
  5 = @(7, 7)
  6 = @(8, 8)
  7 = @(8, 8)
  8 = @(6, 7)
  9 = @(7, 8)
  10 = @(7, 8)
  11 = @(8, 8)
  12 = @(1, 3)
  13 = @(6, 6)
  14 = @(9, 9)
  15 = @(7, 7)
  16 = @(5, 6)
  17 = @(8, 9)
  18 = @(6, 6)
  19 = @(8, 9)
  20 = @(7, 7)
  21 = @(4, 5)
  22 = @(8, 8)
  23 = @(9, 9)
  24 = @(8, 9)
  25 = @(7, 8)
  26 = @(6, 7)
  27 = @(4, 4)
  28 = @(9, 9)
  29 = @(11, 11)
  30 = @(5, 6)
  31 = @(1, 2)
  32 = @(9, 9)
  33 = @(8, 8)
  34 = @(8, 8)
  35 = @(6, 6)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  $ws.Cells.Item($r, 9).Value = $vals[0]
  $ws.Cells.Item($r, 10).Value = $vals[1]
}
